$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "-"
$ws.Range("C3").Value = "-"
$ws.Range("E3").Value = "[-, -, -, 'MCT-3A-Manut. Mecânica']"
$ws.Range("C4").Value = "-"
$ws.Range("E4").Value = "[-, -, -, 'MCT-3A-Manut. Mecânica']"
$ws.Range("E6").Value = "[-, -, -, 'MCT-3A-Manut. Mecânica']"
$ws.Range("E7").Value = "[-, -, -, Carlos Eduardo-Processos de Usinagem 1-3A]"
$ws.Range("C8").Value = "-"

$ws.Range("B10").Value = "-"
$ws.Range("D10").Value = "-"
$ws.Range("B11").Value = "-"
$ws.Range("F11").Value = "['MEC-3A-Retífica', -, -, 'MEC-3A-Mec. Manut.Equip. ind.']"
$ws.Range("B12").Value = "-"
$ws.Range("F12").Value = "['MEC-3A-Retífica', -, -, 'MEC-3A-Mec. Manut.Equip. ind.']"
$ws.Range("B14").Value = "-"
$ws.Range("F14").Value = "['MEC-3A-Retífica', -, -, 'MEC-3A-Mec. Manut.Equip. ind.']"
$ws.Range("B15").Value = "-"
$ws.Range("F15").Value = "['MEC-3A-Retífica', -, -, 'MEC-3A-Mec. Manut.Equip. ind.']"
$ws.Range("B16").Value = "-"
$ws.Range("D16").Value = "-"

$ws.Range("B18").Value = "[-, 'MEC-2NA-M.Maq.E.I.', -, -]"
$ws.Range("C18").Value = "-"
$ws.Range("D18").Value = "-"
$ws.Range("E18").Value = "[-, 'MEC-2NA-Retífica', -, 'ELM-1NA-Manut. Mecânica']"
$ws.Range("F18").Value = "-"

$ws.Range("B19").Value = "[-, 'MEC-2NA-M.Maq.E.I.', -, -]"
$ws.Range("D19").Value = "-"
$ws.Range("E19").Value = "[-, 'MEC-2NA-Retífica', -, 'ELM-1NA-Manut. Mecânica']"

$ws.Range("B20").Value = "[-, 'MEC-2NA-M.Maq.E.I.', -, -]"
$ws.Range("E20").Value = "[-, 'MEC-2NA-Retífica', -, 'ELM-1NA-Manut. Mecânica']"
$ws.Range("F20").Value = "-"

$ws.Range("B21").Value = "[-, 'MEC-2NA-M.Maq.E.I.', -, -]"
$ws.Range("E21").Value = "[-, 'MEC-2NA-Retífica', -, 'ELM-1NA-Manut. Mecânica']"
$ws.Range("F21").Value = "-"
